$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 3 values
$ws.Range("B3").Value = 1000
$ws.Range("D3").Value = 50

# Add new rows 12-15
$ws.Range("A12").Value = "bife do vazio"
$ws.Range("B12").Value = 1000
$ws.Range("C12").Value = "g"
$ws.Range("D12").Value = 60

$ws.Range("A13").Value = "alho poro"
$ws.Range("B13").Value = 100
$ws.Range("C13").Value = "Un"
$ws.Range("D13").Value = 5

$ws.Range("A14").Value = "cebola roxa"
$ws.Range("B14").Value = 500
$ws.Range("C14").Value = "g"
$ws.Range("D14").Value = 7

$ws.Range("A15").Value = "arroz parboilizado"
$ws.Range("B15").Value = 1000
$ws.Range("C15").Value = "Un"
$ws.Range("D15").Value = 30
